$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2015")
$ws.Columns.Item(1).AutoFit()
Write-Host "done"
